$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value = 41369

$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = 0.0625

$excel.CutCopyMode = $false

$ws.Range("B25").Select()
